# "Fixing scripts for watch list"
# Appends two new test-case rows (TestCase_A24, TestCase_A25) to the
# "Test Cases" sheet, right after the existing last row (24).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Clone row 24's formatting (borders etc.) down into the two new rows so
# the new cells pick up the same look as the rest of the table.
$ws.Range("A24:E24").Copy()
$ws.Range("A25:E25").PasteSpecial(-4122)
$ws.Range("A24:E24").Copy()
$ws.Range("A26:E26").PasteSpecial(-4122)

# Row 25 - TestCase_A24
$ws.Range("A25").Value = "TestCase_A24"
$ws.Range("C25").Value = "Verify that TR account gets locked after 5 consecutive unsuccessful login attempts"
$ws.Range("B25").Value = "OPQA-525"
$ws.Range("D25").Value = "Y"
$ws.Range("E25").Value = "SKIP"

# Row 26 - TestCase_A25
$ws.Range("A26").Value = "TestCase_A25"
$ws.Range("B26").Value = "OPQA-529"
$ws.Range("C26").Value = "Verify that Help link is working properly"
$ws.Range("D26").Value = "Y"
$ws.Range("E26").Value = "PASS"

$ws.Range("D19").Select()
